# CDS Library source query fixed
# Replace the Participants-tab query (B2) with the corrected Cypher query,
# and refresh the sheet's look (font size, wrap, row heights, selection)
# to match the latest authoring pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$participantQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE g.library_source in ['Transcriptomic']
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id
LIMIT 100
'@

# 1) Fix the broken library-source filter query used by the Participants tab
$ws.Range("B2").Value = $participantQuery

# 2) Bump the sheet to a consistent 14pt font (was a mix of default 11pt / 12pt)
$ws.Range("A1:E4").Font.Size = 14
$ws.Range("B5:C5").Font.Size = 14
$ws.Range("C6").Font.Size = 14

# 3) Header row + the Tab/Db/Web-file columns (A, D, E) no longer wrap;
#    the query columns (B, C) keep wrapping so the Cypher stays readable.
$ws.Range("A1:E1").WrapText = $false
$ws.Range("A2").WrapText = $false
$ws.Range("D2:E2").WrapText = $false
$ws.Range("A3").WrapText = $false
$ws.Range("D3:E3").WrapText = $false
$ws.Range("A4").WrapText = $false
$ws.Range("D4:E4").WrapText = $false

$ws.Range("B2:C4").WrapText = $true
$ws.Range("B5:C5").WrapText = $true
$ws.Range("C6").WrapText = $true

# 4) Row heights grew to fit the longer/rewrapped query text at 14pt
$ws.Rows(2).RowHeight = 375
$ws.Rows(3).RowHeight = 281.25
$ws.Rows(4).RowHeight = 262.5

# 5) Selection / scroll position saved with the workbook
$ws.Range("A2").Select()
$ws.Range("C2").Select()
